# Update finanza y fotos orlando robert
#
# Adds new "Gastos" (expenses) entry for referee+water, and seven new
# "Ingreso" (income / contribution) entries, then leaves the UI focused
# on the Ingreso sheet at the last entered cell - matching the order of
# operations the author performed (Gastos row first, then Ingreso rows).

$wb = $excel.ActiveWorkbook

# ---- Gastos: one new row (row 81) ----
$gastos = $wb.Worksheets.Item("Gastos")
$gastos.Activate()

$gastos.Range("A81").Value = "2024-02-18"
$gastos.Range("B81").Value = "Arbitro y agua"
$gastos.Range("C81").Formula = "=150+800"

$gastos.Range("A81").Select()

# ---- Ingreso: seven new rows (rows 642-648) ----
$ingreso = $wb.Worksheets.Item("Ingreso")
$ingreso.Activate()

$ingreso.Range("A642").Value = "2024-02-18"
$ingreso.Range("B642").Value = "Rayder"
$ingreso.Range("C642").Value = 100
$ingreso.Range("D642").Value = "Aporte"

$ingreso.Range("A643").Value = "2024-02-18"
$ingreso.Range("B643").Value = "Invitados"
$ingreso.Range("C643").Value = 100
$ingreso.Range("D643").Value = "Aporte"

$ingreso.Range("A644").Value = "2024-02-18"
$ingreso.Range("B644").Value = "Jeicol"
$ingreso.Range("C644").Value = 200
$ingreso.Range("D644").Value = "Aporte"

$ingreso.Range("A645").Value = "2024-02-18"
$ingreso.Range("B645").Value = "Carlos"
$ingreso.Range("C645").Value = 200
$ingreso.Range("D645").Value = "Aporte"

$ingreso.Range("A646").Value = "2024-02-18"
$ingreso.Range("B646").Value = "Rubio"
$ingreso.Range("C646").Value = 100
$ingreso.Range("D646").Value = "Aporte"

$ingreso.Range("A647").Value = "2024-02-18"
$ingreso.Range("B647").Value = "Jordan"
$ingreso.Range("C647").Value = 400
$ingreso.Range("D647").Value = "Aporte"

$ingreso.Range("A648").Value = "2024-02-12"
$ingreso.Range("B648").Value = "Robert"
$ingreso.Range("C648").Value = 300
$ingreso.Range("D648").Value = "Aporte"

# Leave the workbook focused where the author left it: Ingreso sheet,
# selection on the last-entered cell.
$ingreso.Range("B648").Select()
